$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data to match latest scrape.
# Force text format on each target cell first so numeric-looking
# strings (prices, percentages) are stored as text, matching the
# original inline-string cell type instead of being coerced to numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.202.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.01%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.822.66'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.98%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.19%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '452.70'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +7.60%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.28'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +13.92%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.621'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.30%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.13%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.741'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.32%  '
# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.30%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000318'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -7.97%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.69'
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.36'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.29%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.442.44'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.00%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.84%  '
# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.137'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.39%  '
# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.764.04'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.75%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.03'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.06%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.17'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +8.39%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.275.57'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.59%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '427.94'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.24%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.82'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.54%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.24'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +8.43%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.32'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.47%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.48'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +8.84%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '37.21'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.17%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.19'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +20.18%  '
# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.47'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.42%  '
# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.73'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.86%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '736.23'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.20%  '
# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.79'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +12.06%  '
# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.134'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +11.60%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.73'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.78%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '43.10'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +12.58%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.158'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.20%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.27'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.50%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.57'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +13.15%  '
# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.10%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0475'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.67%  '
# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.347'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +13.49%  '
# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.91'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.31%  '
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.92%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.58'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +13.73%  '
# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.17%  '
# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.140'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.81%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.45'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.02%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.26'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +6.54%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.13'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.63%  '
# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.16%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '144.28'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.73%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.87'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.61%  '
